$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1985.375
$ws.Range("I19").Value = 1997.8
$ws.Range("J19").Value = 1964.6666
$ws.Range("K19").Value = 1997.8
$ws.Range("L19").Value = 1964.6666
$ws.Range("M19").Value = -1822.8
$ws.Range("N19").Value = -2314.6666
$ws.Range("H94").Value = 20413008
$ws.Range("I94").Value = 20413008
$ws.Range("K94").Value = 20413008
$ws.Range("M94").Value = -20412557
$ws.Range("H98").Value = 2097.3684
$ws.Range("I98").Value = 1803
$ws.Range("K98").Value = 1803
$ws.Range("M98").Value = -305
$ws.Range("H122").Value = 2097.3684
$ws.Range("I122").Value = 1803
$ws.Range("K122").Value = 5409
$ws.Range("M122").Value = -2959
$ws.Range("H137").Value = 10377.583
$ws.Range("I137").Value = 16994.5
$ws.Range("J137").Value = 3760.6667
$ws.Range("K137").Value = 50983.5
$ws.Range("L137").Value = 11282.0001
$ws.Range("M137").Value = -48433.5
$ws.Range("N137").Value = -16382.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7156.3804
$ws.Range("I32").Value = 6834.403
$ws.Range("K32").Value = 6834.403
$ws.Range("M32").Value = -6547.403
$ws.Range("H39").Value = 2697.8
$ws.Range("J39").Value = 2450
$ws.Range("L39").Value = 2450
$ws.Range("N39").Value = -3490
$ws.Range("H61").Value = 4417.375
$ws.Range("I61").Value = 4681.5835
$ws.Range("K61").Value = 4681.5835
$ws.Range("M61").Value = -4469.5835
$ws.Range("H102").Value = 4435.7896
$ws.Range("I102").Value = 2479.75
$ws.Range("K102").Value = 2479.75
$ws.Range("M102").Value = -857.75
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N111").Value = 0
$ws.Range("H132").Value = 2452.9023
$ws.Range("I132").Value = 2384.7646
$ws.Range("K132").Value = 7154.293799999999
$ws.Range("M132").Value = -4624.293799999999
$ws.Range("H133").Value = 81856.57000000001
$ws.Range("J133").Value = 105880.5
$ws.Range("L133").Value = 105880.5
$ws.Range("N133").Value = -110940.5
$ws.Range("H136").Value = 4417.375
$ws.Range("I136").Value = 4681.5835
$ws.Range("K136").Value = 14044.7505
$ws.Range("M136").Value = -11494.7505
$ws.Range("L111").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 99998
$ws.Range("J13").Value = 99998
$ws.Range("L13").Value = 99998
$ws.Range("N13").Value = -100334
$ws.Range("H80").Value = 5295.4546
$ws.Range("I80").Value = 13328.25
$ws.Range("J80").Value = 705.2857
$ws.Range("K80").Value = 13328.25
$ws.Range("L80").Value = 705.2857
$ws.Range("M80").Value = -12330.25
$ws.Range("N80").Value = -2701.2857
$ws.Range("H83").Value = 5295.4546
$ws.Range("I83").Value = 13328.25
$ws.Range("J83").Value = 705.2857
$ws.Range("K83").Value = 66641.25
$ws.Range("L83").Value = 3526.4285
$ws.Range("M83").Value = -61649.25
$ws.Range("N83").Value = -13510.4285
$ws.Range("H135").Value = 192713.14
$ws.Range("J135").Value = 192713.14
$ws.Range("L135").Value = 192713.14
$ws.Range("N135").Value = -202853.14

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3185.4546
$ws.Range("I31").Value = 2580.125
$ws.Range("K31").Value = 2580.125
$ws.Range("M31").Value = -2285.125
$ws.Range("H34").Value = 3185.4546
$ws.Range("I34").Value = 2580.125
$ws.Range("K34").Value = 2580.125
$ws.Range("M34").Value = -2378.125
$ws.Range("H86").Value = 9527552
$ws.Range("I86").Value = 11114644
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 11114644
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -11113521
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 9527552
$ws.Range("I89").Value = 11114644
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 55573220
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -55567604
$ws.Range("N89").Value = -36232
$ws.Range("H105").Value = 1353.2222
$ws.Range("I105").Value = 1296.4667
$ws.Range("K105").Value = 1296.4667
$ws.Range("M105").Value = 450.5333000000001
$ws.Range("H132").Value = 5790.9346
$ws.Range("I132").Value = 1834.4546
$ws.Range("J132").Value = 15834.308
$ws.Range("K132").Value = 5503.3638
$ws.Range("L132").Value = 47502.924
$ws.Range("M132").Value = -2973.3638
$ws.Range("N132").Value = -52562.924
$ws.Range("H138").Value = 69502.914
$ws.Range("J138").Value = 71275.91
$ws.Range("L138").Value = 71275.91
$ws.Range("N138").Value = -81555.91
$ws.Range("H141").Value = 153512.88
$ws.Range("J141").Value = 172381.9
$ws.Range("L141").Value = 172381.9
$ws.Range("N141").Value = -182741.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 14501.667
$ws.Range("J82").Value = 14104.667
$ws.Range("L82").Value = 42314.001
$ws.Range("N82").Value = -43126.001
$ws.Range("H85").Value = 14501.667
$ws.Range("J85").Value = 14104.667
$ws.Range("L85").Value = 42314.001
$ws.Range("N85").Value = -45122.001
$ws.Range("H103").Value = 531.2
$ws.Range("I103").Value = 451.66666
$ws.Range("K103").Value = 1354.99998
$ws.Range("M103").Value = -475.9999800000001
$ws.Range("H128").Value = 282165.84
$ws.Range("I128").Value = 282165.84
$ws.Range("K128").Value = 846497.52
$ws.Range("M128").Value = -841517.52
$ws.Range("H132").Value = 2708.611
$ws.Range("I132").Value = 1910.25
$ws.Range("K132").Value = 17192.25
$ws.Range("M132").Value = -14662.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 120002000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 120002000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H97").Value = 1164.5
$ws.Range("I97").Value = 723.8
$ws.Range("K97").Value = 723.8
$ws.Range("M97").Value = -227.8
$ws.Range("H122").Value = 4012.182
$ws.Range("I122").Value = 3514.889
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 10544.667
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -8094.667000000001
$ws.Range("N122").Value = -23650
$ws.Range("H136").Value = 25582.785
$ws.Range("J136").Value = 25582.785
$ws.Range("L136").Value = 76748.355
$ws.Range("N136").Value = -81848.355
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4990
$ws.Range("I7").Value = 4990
$ws.Range("K7").Value = 4990
$ws.Range("M7").Value = -4878
$ws.Range("H20").Value = 9070
$ws.Range("J20").Value = 9775
$ws.Range("L20").Value = 9775
$ws.Range("N20").Value = -10227
$ws.Range("H40").Value = 4564.9585
$ws.Range("I40").Value = 2913.9285
$ws.Range("J40").Value = 6876.4
$ws.Range("K40").Value = 2913.9285
$ws.Range("L40").Value = 6876.4
$ws.Range("M40").Value = -2777.9285
$ws.Range("N40").Value = -7148.4
$ws.Range("H46").Value = 2949.5625
$ws.Range("J46").Value = 3763
$ws.Range("L46").Value = 3763
$ws.Range("N46").Value = -4139
$ws.Range("H126").Value = 4990
$ws.Range("I126").Value = 4990
$ws.Range("K126").Value = 14970
$ws.Range("M126").Value = -12500
$ws.Range("H132").Value = 3029.0286
$ws.Range("I132").Value = 2695.7083
$ws.Range("K132").Value = 8087.124899999999
$ws.Range("M132").Value = -5557.124899999999
$ws.Range("H135").Value = 78448
$ws.Range("J135").Value = 78448
$ws.Range("L135").Value = 78448
$ws.Range("N135").Value = -88588
$ws.Range("H137").Value = 72523
$ws.Range("J137").Value = 72523
$ws.Range("L137").Value = 72523
$ws.Range("N137").Value = -82723

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 49626
$ws.Range("I49").Value = 49626
$ws.Range("K49").Value = 49626
$ws.Range("M49").Value = -49396
$ws.Range("H118").Value = 109999
$ws.Range("J118").Value = 109999
$ws.Range("L118").Value = 109999
$ws.Range("N118").Value = -113313
